# edit.ps1 — applies the 2023-11-12 14:45 betexplorer scrape refresh to Sheet1
# Source match rows were re-sequenced within several match-day blocks (fixture
# order changed upstream) and three fixtures played on 2023-11-11 were appended
# (rows 134-136), extending the sheet from 133 to 136 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: re-synchronise the home/away fixture data (columns F:V) for the rows
# whose content was swapped/rotated between sibling rows of the same match day.
# (columns A:E -- index/country/tournament/season/date -- are untouched)
# ---------------------------------------------------------------------------
# row 5  (was row 6: Bilje vs Ilirija)
$ws.Range("F5").Value = "Bilje"
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = "Ilirija"
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 3.49
$ws.Range("K5").Value = "28/07/2023 06:12"
$ws.Range("L5").Value = 2.05
$ws.Range("M5").Value = "29/07/2023 17:59"
$ws.Range("N5").Value = 3.43
$ws.Range("O5").Value = "28/07/2023 06:12"
$ws.Range("P5").Value = 3.62
$ws.Range("Q5").Value = "29/07/2023 17:59"
$ws.Range("R5").Value = 1.85
$ws.Range("S5").Value = "28/07/2023 06:12"
$ws.Range("T5").Value = 2.72
$ws.Range("U5").Value = "29/07/2023 17:59"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-ilirija/6PRxXGfF/"

# row 6  (was row 5: Grosuplje vs ND Gorica)
$ws.Range("F6").Value = "Grosuplje"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "ND Gorica"
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 3.35
$ws.Range("K6").Value = "29/07/2023 14:12"
$ws.Range("L6").Value = 3.5
$ws.Range("M6").Value = "29/07/2023 17:02"
$ws.Range("N6").Value = 3.29
$ws.Range("O6").Value = "29/07/2023 14:12"
$ws.Range("P6").Value = 3.16
$ws.Range("Q6").Value = "29/07/2023 17:02"
$ws.Range("R6").Value = 2.04
$ws.Range("S6").Value = "29/07/2023 14:12"
$ws.Range("T6").Value = 2.07
$ws.Range("U6").Value = "29/07/2023 17:02"
$ws.Range("V6").Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-nd-gorica/O0xYXdu9/"

# row 14  (was row 15: Primorje vs Bilje)
$ws.Range("F14").Value = "Primorje"
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = "Bilje"
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1.64
$ws.Range("K14").Value = "12/08/2023 05:43"
$ws.Range("L14").Value = 1.94
$ws.Range("M14").Value = "13/08/2023 16:02"
$ws.Range("N14").Value = 3.65
$ws.Range("O14").Value = "12/08/2023 05:43"
$ws.Range("P14").Value = 3.51
$ws.Range("Q14").Value = "13/08/2023 16:02"
$ws.Range("R14").Value = 4.07
$ws.Range("S14").Value = "12/08/2023 05:43"
$ws.Range("T14").Value = 3.5
$ws.Range("U14").Value = "13/08/2023 16:02"
$ws.Range("V14").Value = "https://www.betexplorer.com/football/slovenia/2-snl/primorje-bilje/b3Awl9V6/"

# row 15  (was row 14: ND Gorica vs Nafta)
$ws.Range("F15").Value = "ND Gorica"
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = "Nafta"
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 2.03
$ws.Range("K15").Value = "13/08/2023 10:40"
$ws.Range("L15").Value = 2.25
$ws.Range("M15").Value = "13/08/2023 15:23"
$ws.Range("N15").Value = 3.49
$ws.Range("O15").Value = "13/08/2023 10:40"
$ws.Range("P15").Value = 3.42
$ws.Range("Q15").Value = "13/08/2023 15:34"
$ws.Range("R15").Value = 3.21
$ws.Range("S15").Value = "13/08/2023 10:40"
$ws.Range("T15").Value = 2.83
$ws.Range("U15").Value = "13/08/2023 15:23"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nd-gorica-nafta/xlMUkV0f/"

# row 24  (was row 25: Dravinja vs Ilirija)
$ws.Range("F24").Value = "Dravinja"
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = "Ilirija"
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 3.12
$ws.Range("K24").Value = "20/08/2023 08:59"
$ws.Range("L24").Value = 3.08
$ws.Range("M24").Value = "20/08/2023 11:20"
$ws.Range("N24").Value = 3.27
$ws.Range("O24").Value = "20/08/2023 08:59"
$ws.Range("P24").Value = 3.4
$ws.Range("Q24").Value = "20/08/2023 15:35"
$ws.Range("R24").Value = 2.15
$ws.Range("S24").Value = "20/08/2023 08:59"
$ws.Range("T24").Value = 2.12
$ws.Range("U24").Value = "20/08/2023 14:26"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/slovenia/2-snl/dravinja-ilirija/zRWDoyR5/"

# row 25  (was row 24: Bilje vs Tabor Sezana)
$ws.Range("F25").Value = "Bilje"
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = "Tabor Sezana"
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 1.52
$ws.Range("K25").Value = "20/08/2023 09:00"
$ws.Range("L25").Value = 1.47
$ws.Range("M25").Value = "20/08/2023 17:19"
$ws.Range("N25").Value = 4.26
$ws.Range("O25").Value = "20/08/2023 09:00"
$ws.Range("P25").Value = 5.59
$ws.Range("Q25").Value = "20/08/2023 17:19"
$ws.Range("R25").Value = 5.03
$ws.Range("S25").Value = "20/08/2023 09:00"
$ws.Range("T25").Value = 4.39
$ws.Range("U25").Value = "20/08/2023 17:25"
$ws.Range("V25").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-tabor-sezana/6ouckJRt/"

# row 69  (was row 70: Triglav vs ND Gorica)
$ws.Range("F69").Value = "Triglav"
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = "ND Gorica"
$ws.Range("I69").Value = 3
$ws.Range("J69").Value = 3.01
$ws.Range("K69").Value = "23/09/2023 03:12"
$ws.Range("L69").Value = 3.77
$ws.Range("M69").Value = "24/09/2023 15:32"
$ws.Range("N69").Value = 3.29
$ws.Range("O69").Value = "23/09/2023 03:12"
$ws.Range("P69").Value = 3.27
$ws.Range("Q69").Value = "24/09/2023 15:39"
$ws.Range("R69").Value = 2.03
$ws.Range("S69").Value = "23/09/2023 03:12"
$ws.Range("T69").Value = 1.94
$ws.Range("U69").Value = "24/09/2023 15:39"
$ws.Range("V69").Value = "https://www.betexplorer.com/football/slovenia/2-snl/triglav-nd-gorica/IsAXSPF2/"

# row 70  (was row 69: Tolmin vs Jadran Dekani)
$ws.Range("F70").Value = "Tolmin"
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = "Jadran Dekani"
$ws.Range("I70").Value = 2
$ws.Range("J70").Value = 3.15
$ws.Range("K70").Value = "23/09/2023 03:12"
$ws.Range("L70").Value = 3.52
$ws.Range("M70").Value = "24/09/2023 15:42"
$ws.Range("N70").Value = 3.16
$ws.Range("O70").Value = "23/09/2023 03:12"
$ws.Range("P70").Value = 3.17
$ws.Range("Q70").Value = "24/09/2023 15:41"
$ws.Range("R70").Value = 2.06
$ws.Range("S70").Value = "23/09/2023 03:12"
$ws.Range("T70").Value = 2.06
$ws.Range("U70").Value = "24/09/2023 15:41"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/slovenia/2-snl/tolmin-jadran-dekani/GxhKJLxk/"

# row 73  (was row 74: ND Gorica vs Fuzinar)
$ws.Range("F73").Value = "ND Gorica"
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = "Fuzinar"
$ws.Range("I73").Value = 1
$ws.Range("J73").Value = 1.29
$ws.Range("K73").Value = "28/09/2023 02:42"
$ws.Range("L73").Value = 1.25
$ws.Range("M73").Value = "29/09/2023 13:45"
$ws.Range("N73").Value = 4.98
$ws.Range("O73").Value = "28/09/2023 02:42"
$ws.Range("P73").Value = 6.11
$ws.Range("Q73").Value = "29/09/2023 15:29"
$ws.Range("R73").Value = 6.88
$ws.Range("S73").Value = "28/09/2023 02:42"
$ws.Range("T73").Value = 7.81
$ws.Range("U73").Value = "29/09/2023 15:29"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nd-gorica-fuzinar/xQJeJCr6/"

# row 74  (was row 73: Ilirija vs Primorje)
$ws.Range("F74").Value = "Ilirija"
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = "Primorje"
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3.72
$ws.Range("K74").Value = "28/09/2023 02:42"
$ws.Range("L74").Value = 4.94
$ws.Range("M74").Value = "29/09/2023 15:20"
$ws.Range("N74").Value = 3.4
$ws.Range("O74").Value = "28/09/2023 02:42"
$ws.Range("P74").Value = 3.8
$ws.Range("Q74").Value = "29/09/2023 15:20"
$ws.Range("R74").Value = 1.79
$ws.Range("S74").Value = "28/09/2023 02:42"
$ws.Range("T74").Value = 1.61
$ws.Range("U74").Value = "29/09/2023 15:20"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-primorje/f1UjKhTa/"

# row 76  (was row 79: Rudar vs Triglav)
$ws.Range("F76").Value = "Rudar"
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = "Triglav"
$ws.Range("I76").Value = 2
$ws.Range("J76").Value = 2.29
$ws.Range("K76").Value = "29/09/2023 02:42"
$ws.Range("L76").Value = 2.25
$ws.Range("M76").Value = "30/09/2023 15:29"
$ws.Range("N76").Value = 3.2
$ws.Range("O76").Value = "29/09/2023 02:42"
$ws.Range("P76").Value = 3.37
$ws.Range("Q76").Value = "30/09/2023 15:25"
$ws.Range("R76").Value = 2.64
$ws.Range("S76").Value = "29/09/2023 02:42"
$ws.Range("T76").Value = 2.9
$ws.Range("U76").Value = "30/09/2023 15:24"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/slovenia/2-snl/rudar-triglav/jepXGJM7/"

# row 78  (was row 76: NK Bistrica vs Tabor Sezana)
$ws.Range("F78").Value = "NK Bistrica"
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = "Tabor Sezana"
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1.41
$ws.Range("K78").Value = "29/09/2023 02:42"
$ws.Range("L78").Value = 1.34
$ws.Range("M78").Value = "30/09/2023 15:24"
$ws.Range("N78").Value = 4.29
$ws.Range("O78").Value = "29/09/2023 02:42"
$ws.Range("P78").Value = 5.2
$ws.Range("Q78").Value = "30/09/2023 15:29"
$ws.Range("R78").Value = 5.33
$ws.Range("S78").Value = "29/09/2023 02:42"
$ws.Range("T78").Value = 6.9
$ws.Range("U78").Value = "30/09/2023 15:29"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bistrica-tabor-sezana/tbQnLYDg/"

# row 79  (was row 78: Nafta vs Bilje)
$ws.Range("F79").Value = "Nafta"
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = "Bilje"
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = 1.72
$ws.Range("K79").Value = "29/09/2023 02:42"
$ws.Range("L79").Value = 1.78
$ws.Range("M79").Value = "30/09/2023 15:23"
$ws.Range("N79").Value = 3.77
$ws.Range("O79").Value = "29/09/2023 02:42"
$ws.Range("P79").Value = 4.06
$ws.Range("Q79").Value = "30/09/2023 15:23"
$ws.Range("R79").Value = 3.65
$ws.Range("S79").Value = "29/09/2023 02:42"
$ws.Range("T79").Value = 3.61
$ws.Range("U79").Value = "30/09/2023 15:23"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nafta-bilje/4bYEu3a8/"

# row 98  (was row 100: Tabor Sezana vs Fuzinar)
$ws.Range("F98").Value = "Tabor Sezana"
$ws.Range("G98").Value = 4
$ws.Range("H98").Value = "Fuzinar"
$ws.Range("I98").Value = 4
$ws.Range("J98").Value = 2.71
$ws.Range("K98").Value = "19/10/2023 02:12"
$ws.Range("L98").Value = 3.21
$ws.Range("M98").Value = "20/10/2023 14:55"
$ws.Range("N98").Value = 3.31
$ws.Range("O98").Value = "19/10/2023 02:12"
$ws.Range("P98").Value = 3.81
$ws.Range("Q98").Value = "20/10/2023 14:59"
$ws.Range("R98").Value = 2.19
$ws.Range("S98").Value = "19/10/2023 02:12"
$ws.Range("T98").Value = 1.95
$ws.Range("U98").Value = "20/10/2023 14:59"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/slovenia/2-snl/tabor-sezana-fuzinar/I33Cpks5/"

# row 100  (was row 98: Grosuplje vs NK Bistrica)
$ws.Range("F100").Value = "Grosuplje"
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = "NK Bistrica"
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = 1.75
$ws.Range("K100").Value = "19/10/2023 02:12"
$ws.Range("L100").Value = 1.75
$ws.Range("M100").Value = "20/10/2023 14:57"
$ws.Range("N100").Value = 3.44
$ws.Range("O100").Value = "19/10/2023 02:12"
$ws.Range("P100").Value = 3.56
$ws.Range("Q100").Value = "20/10/2023 14:57"
$ws.Range("R100").Value = 3.74
$ws.Range("S100").Value = "19/10/2023 02:12"
$ws.Range("T100").Value = 4.32
$ws.Range("U100").Value = "20/10/2023 14:57"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-bistrica/hMePsmRN/"

# row 117  (was row 118: ND Gorica vs Tolmin)
$ws.Range("F117").Value = "ND Gorica"
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = "Tolmin"
$ws.Range("I117").Value = 3
$ws.Range("J117").Value = 1.27
$ws.Range("K117").Value = "28/10/2023 03:12"
$ws.Range("L117").Value = 1.32
$ws.Range("M117").Value = "29/10/2023 11:36"
$ws.Range("N117").Value = 5
$ws.Range("O117").Value = "28/10/2023 03:12"
$ws.Range("P117").Value = 4.94
$ws.Range("Q117").Value = "29/10/2023 12:03"
$ws.Range("R117").Value = 7.14
$ws.Range("S117").Value = "28/10/2023 03:12"
$ws.Range("T117").Value = 7.46
$ws.Range("U117").Value = "29/10/2023 11:36"
$ws.Range("V117").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nd-gorica-tolmin/ncz2e6eN/"

# row 118  (was row 117: NK Krka vs Rudar)
$ws.Range("F118").Value = "NK Krka"
$ws.Range("G118").Value = 4
$ws.Range("H118").Value = "Rudar"
$ws.Range("I118").Value = 5
$ws.Range("J118").Value = 1.9
$ws.Range("K118").Value = "28/10/2023 03:12"
$ws.Range("L118").Value = 2.06
$ws.Range("M118").Value = "29/10/2023 13:33"
$ws.Range("N118").Value = 3.4
$ws.Range("O118").Value = "28/10/2023 03:12"
$ws.Range("P118").Value = 3.61
$ws.Range("Q118").Value = "29/10/2023 13:33"
$ws.Range("R118").Value = 3.25
$ws.Range("S118").Value = "28/10/2023 03:12"
$ws.Range("T118").Value = 3.09
$ws.Range("U118").Value = "29/10/2023 13:33"
$ws.Range("V118").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nk-krka-rudar/d8skb8B4/"

# row 130  (was row 132: Bilje vs ND Gorica)
$ws.Range("F130").Value = "Bilje"
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = "ND Gorica"
$ws.Range("I130").Value = 1
$ws.Range("J130").Value = 3.31
$ws.Range("K130").Value = "10/11/2023 02:13"
$ws.Range("L130").Value = 4.22
$ws.Range("M130").Value = "11/11/2023 13:59"
$ws.Range("N130").Value = 3.41
$ws.Range("O130").Value = "10/11/2023 02:13"
$ws.Range("P130").Value = 4.12
$ws.Range("Q130").Value = "11/11/2023 13:59"
$ws.Range("R130").Value = 1.87
$ws.Range("S130").Value = "10/11/2023 02:13"
$ws.Range("T130").Value = 1.65
$ws.Range("U130").Value = "11/11/2023 13:59"
$ws.Range("V130").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-nd-gorica/tKHAlwfE/"

# row 131  (was row 130: Dravinja vs Triglav)
$ws.Range("F131").Value = "Dravinja"
$ws.Range("G131").Value = 2
$ws.Range("H131").Value = "Triglav"
$ws.Range("I131").Value = 1
$ws.Range("J131").Value = 2.58
$ws.Range("K131").Value = "10/11/2023 02:13"
$ws.Range("L131").Value = 2.8
$ws.Range("M131").Value = "11/11/2023 13:39"
$ws.Range("N131").Value = 3.09
$ws.Range("O131").Value = "10/11/2023 02:13"
$ws.Range("P131").Value = 3.34
$ws.Range("Q131").Value = "11/11/2023 13:39"
$ws.Range("R131").Value = 2.4
$ws.Range("S131").Value = "10/11/2023 02:13"
$ws.Range("T131").Value = 2.34
$ws.Range("U131").Value = "11/11/2023 13:39"
$ws.Range("V131").Value = "https://www.betexplorer.com/football/slovenia/2-snl/dravinja-triglav/dSw6jH9e/"

# row 132  (was row 131: Grosuplje vs Fuzinar)
$ws.Range("F132").Value = "Grosuplje"
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = "Fuzinar"
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1.39
$ws.Range("K132").Value = "10/11/2023 02:13"
$ws.Range("L132").Value = 1.56
$ws.Range("M132").Value = "11/11/2023 13:56"
$ws.Range("N132").Value = 4.33
$ws.Range("O132").Value = "10/11/2023 02:13"
$ws.Range("P132").Value = 4.28
$ws.Range("Q132").Value = "11/11/2023 13:57"
$ws.Range("R132").Value = 5.51
$ws.Range("S132").Value = "10/11/2023 02:13"
$ws.Range("T132").Value = 4.69
$ws.Range("U132").Value = "11/11/2023 13:57"
$ws.Range("V132").Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-fuzinar/zsHEmc9K/"

# ---------------------------------------------------------------------------
# Part 2: append the three new fixtures (rows 134-136) scraped for 2023-11-11
# ---------------------------------------------------------------------------
foreach ($r in 134..136) {
    $ws.Range("A133").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("E133").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)
}

# row 134  (Rudar 2 - 1 Ilirija)
$ws.Range("A134").Value = 133
$ws.Range("B134").Value = "slovenia"
$ws.Range("C134").Value = "2-snl"
$ws.Range("D134").Value = "2023-2024"
$ws.Range("E134").Value = 45242.58333333334
$ws.Range("F134").Value = "Rudar"
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = "Ilirija"
$ws.Range("I134").Value = 1
$ws.Range("J134").Value = 1.85
$ws.Range("K134").Value = "11/11/2023 02:13"
$ws.Range("L134").Value = 2.02
$ws.Range("M134").Value = "12/11/2023 13:57"
$ws.Range("N134").Value = 3.44
$ws.Range("O134").Value = "11/11/2023 02:13"
$ws.Range("P134").Value = 3.86
$ws.Range("Q134").Value = "12/11/2023 13:57"
$ws.Range("R134").Value = 3.46
$ws.Range("S134").Value = "11/11/2023 02:13"
$ws.Range("T134").Value = 3.01
$ws.Range("U134").Value = "12/11/2023 13:57"
$ws.Range("V134").Value = "https://www.betexplorer.com/football/slovenia/2-snl/rudar-ilirija/2TI6kJv8/"

# row 135  (Beltinci 2 - 0 Tabor Sezana)
$ws.Range("A135").Value = 134
$ws.Range("B135").Value = "slovenia"
$ws.Range("C135").Value = "2-snl"
$ws.Range("D135").Value = "2023-2024"
$ws.Range("E135").Value = 45242.58333333334
$ws.Range("F135").Value = "Beltinci"
$ws.Range("G135").Value = 2
$ws.Range("H135").Value = "Tabor Sezana"
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1.2
$ws.Range("K135").Value = "11/11/2023 02:13"
$ws.Range("L135").Value = 1.23
$ws.Range("M135").Value = "12/11/2023 13:57"
$ws.Range("N135").Value = 5.72
$ws.Range("O135").Value = "11/11/2023 02:13"
$ws.Range("P135").Value = 6.45
$ws.Range("Q135").Value = "12/11/2023 13:58"
$ws.Range("R135").Value = 8.22
$ws.Range("S135").Value = "11/11/2023 02:13"
$ws.Range("T135").Value = 8.2
$ws.Range("U135").Value = "12/11/2023 13:58"
$ws.Range("V135").Value = "https://www.betexplorer.com/football/slovenia/2-snl/beltinci-tabor-sezana/fRcoryXs/"

# row 136  (Tolmin 0 - 0 Primorje)
$ws.Range("A136").Value = 135
$ws.Range("B136").Value = "slovenia"
$ws.Range("C136").Value = "2-snl"
$ws.Range("D136").Value = "2023-2024"
$ws.Range("E136").Value = 45242.58333333334
$ws.Range("F136").Value = "Tolmin"
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = "Primorje"
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3.95
$ws.Range("K136").Value = "11/11/2023 02:13"
$ws.Range("L136").Value = 4.6
$ws.Range("M136").Value = "12/11/2023 13:51"
$ws.Range("N136").Value = 3.54
$ws.Range("O136").Value = "11/11/2023 02:13"
$ws.Range("P136").Value = 3.83
$ws.Range("Q136").Value = "12/11/2023 13:51"
$ws.Range("R136").Value = 1.71
$ws.Range("S136").Value = "11/11/2023 02:13"
$ws.Range("T136").Value = 1.65
$ws.Range("U136").Value = "12/11/2023 13:51"
$ws.Range("V136").Value = "https://www.betexplorer.com/football/slovenia/2-snl/tolmin-primorje/lpLInHOQ/"

